$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark that currently sits after "Frontend" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Locate the "Karten design, wie?" paragraph and rebuild it ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $txt = $cand.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Karten design, wie?") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$full = $d.Range($target.Range.Start, $target.Range.End)
$full.Delete()

$afterDeleteParagraph = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $d.Range($afterDeleteParagraph.Range.Start, $afterDeleteParagraph.Range.Start)
$insertionPoint.InsertParagraphBefore()

$newParagraph = $d.Paragraphs.Item($targetIndex)
$newRange = $d.Range($newParagraph.Range.Start, $newParagraph.Range.End)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="000661EE" w:rsidRPr="00274828" w:rsidRDefault="000661EE" w:rsidP="000661EE">
<w:pPr>
<w:pStyle w:val="Listenabsatz"/>
<w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr>
<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:u w:val="single"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Karten </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>design</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>, wie?</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> (https://developers.google.com/maps/documentation/javascript/overview)</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newRange.InsertXML($xmlFrag)
